$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet to reflect the new "through" date.
$ws.Name = "Through 2022-05-07"

# Update the row label for the May data row (row 6).
$ws.Range("A6").Value = "May (through 05-07)"

# April row (row 5) - only the 2022 column (I) changes.
$ws.Range("I5").Value = 117

# May row (row 6) - updated counts for several years.
$ws.Range("B6").Value = 5
$ws.Range("D6").Value = 12
$ws.Range("H6").Value = 27
$ws.Range("I6").Value = 24

# Total row (row 7) - updated sums.
$ws.Range("B7").Value = 94
$ws.Range("D7").Value = 265
$ws.Range("H7").Value = 550
$ws.Range("I7").Value = 576
